# Crypto price/volume data refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "62.722.47"
    "E2" = "  +3.16%  "
    "D3" = "2.445.30"
    "E3" = "  +1.89%  "
    "E4" = "  -0.14%  "
    "D5" = "578.02"
    "E5" = "  +2.95%  "
    "D6" = "145.87"
    "E6" = "  +2.78%  "
    "E7" = "  +0.07%  "
    "E8" = "  +0.15%  "
    "D9" = "2.444.09"
    "E9" = "  +1.56%  "
    "E10" = "  +2.76%  "
    "E11" = "  +2.03%  "
    "E12" = "  +0.82%  "
    "E13" = "  +2.37%  "
    "D14" = "28.21"
    "E14" = "  +7.95%  "
    "E15" = "  +5.51%  "
    "E16" = "  +1.65%  "
    "D17" = "62.637.92"
    "E17" = "  +3.31%  "
    "D18" = "2.452.15"
    "E18" = "  +1.81%  "
    "D19" = "0.0₆0943"
    "E19" = "  +240.19%  "
    "E20" = "  -3.96%  "
    "D21" = "10.94"
    "E21" = "  +2.51%  "
    "D22" = "329.68"
    "E22" = "  +1.83%  "
    "E23" = "  +1.12%  "
    "E24" = "  +10.02%  "
    "E25" = "  +0.00%  "
    "D26" = "65.84"
    "E26" = "  +1.87%  "
    "D27" = "639.76"
    "E27" = "  +11.56%  "
    "D28" = "1.18"
    "E28" = "  +17.64%  "
    "D29" = "8.45"
    "E29" = "  +5.27%  "
    "D30" = "0.0₃0984"
    "E30" = "  +4.83%  "
    "D31" = "2.565.81"
    "B32" = "Fetch.AI"
    "C32" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "D32" = "1.45"
    "E32" = "  +8.33%  "
    "B33" = "InternetComputer(DFINITY)"
    "C33" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D33" = "8.19"
    "E33" = "  +1.49%  "
    "D34" = "1.87"
    "E34" = "  +3.51%  "
    "E35" = "  +4.55%  "
    "E36" = "  +2.45%  "
    "E37" = "  +0.14%  "
    "E38" = "  +3.32%  "
    "D39" = "5.50"
    "E39" = "  +6.53%  "
    "E40" = "  +0.86%  "
    "D41" = "152.12"
    "E41" = "  -0.03%  "
    "D42" = "18.70"
    "E42" = "  +2.27%  "
    "E43" = "  +5.64%  "
    "D44" = "1.76"
    "D45" = "42.26"
    "E45" = "  +1.28%  "
    "D47" = "14.95"
    "E47" = "  +27.52%  "
    "D48" = "145.24"
    "E48" = "  +2.12%  "
    "E49" = "  +2.41%  "
    "D50" = "20.60"
    "E50" = "  +6.58%  "
    "E51" = "  +2.95%  "
}

# Cells whose new value is a plain number-looking string need to be forced to
# Text format first, otherwise Excel auto-converts them to a Number type and
# they would lose formatting (these are price columns stored as text).
$textForceCells = @(
    "D5"
    "D6"
    "D14"
    "D21"
    "D22"
    "D26"
    "D27"
    "D28"
    "D29"
    "D32"
    "D33"
    "D34"
    "D39"
    "D41"
    "D42"
    "D44"
    "D45"
    "D47"
    "D48"
    "D50"
)

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    if ($textForceCells -contains $cellRef) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$cellRef]
}
